$wb = $excel.ActiveWorkbook

# "Generate Report for Handoff": the handoff run completed for the
# remaining (previously un-handed-off) rows, so their Priority flips
# from "low" to "ht" and their Latest Handoff Datetime is refreshed.

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E4").Value = "ht"
$zhcn.Range("E5").Value = "ht"
$zhcn.Range("E6").Value = "ht"
$zhcn.Range("E7").Value = "ht"

$zhcn.Range("H4").Value = "2016-09-01 04:35:16"
$zhcn.Range("H5").Value = "2016-09-01 04:35:16"
$zhcn.Range("H6").Value = "2016-09-01 04:35:16"
$zhcn.Range("H7").Value = "2016-09-01 04:35:16"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E4").Value = "ht"
$dede.Range("E5").Value = "ht"
$dede.Range("E6").Value = "ht"
$dede.Range("E7").Value = "ht"

$dede.Range("H4").Value = "2016-09-01 04:35:20"
$dede.Range("H5").Value = "2016-09-01 04:35:20"
$dede.Range("H6").Value = "2016-09-01 04:35:20"
$dede.Range("H7").Value = "2016-09-01 04:35:20"

# The Overview sheet mirrors each file's "Latest HO Xliff Generate Date"
# from its de-de handoff timestamp, so it picks up the same refresh.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G4").Value = "2016-09-01 04:35:20"
$overview.Range("G5").Value = "2016-09-01 04:35:20"
$overview.Range("G6").Value = "2016-09-01 04:35:20"
$overview.Range("G7").Value = "2016-09-01 04:35:20"
